# TestData.xlsx — "Added additional test case as a result of a bug" fix-up.
#
# The UserName column (D) on the UserDetails sheet stores long numeric IDs
# but the column is formatted as Text ("@"), so the two existing data rows
# get their numeric literals re-randomised:
#   D2: 5697768474 -> 4876941163
#   D3: 8857887012 -> 5212638910
#
# Because the column's cell style is Text, assigning a bare number through
# .Value would be auto-coerced to a text/shared-string cell (losing the
# numeric <c t="n"> storage the workbook already uses for D2/D3). To keep
# the values stored as genuine numbers - exactly like the original file -
# we briefly flip the number format to a numeric one, write the values,
# then restore the original Text format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Range("D2:D3")
$originalFormat = $col.NumberFormat

$col.NumberFormat = "0"
$ws.Range("D2").Value = 4876941163
$ws.Range("D3").Value = 5212638910
$col.NumberFormat = $originalFormat
